$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5,9,11,17,23,24,28,29,39,50,53,57,66,81,89,99,105,128,140,150,151,156,164,166,167,169,173,194,195,198,199,206,214,216,220,221,225,226,230,236,240,248,249,255,260,261,265,268,269,282,283,292,298,300,305,314,315,330,333,336,358,363,385,387,394,401,403,406,422,432,438,440,441,442,450,454,460,461,462,466,468,471,475,476,485,488,490,494,498,500,506,507,510,513,515,518,519,523,526,528,531,545,549,552,553,554,557,560,593,596,603,611,613,614,619,626,632,637,657,660,671,674,676,681,682,685,687,690,694,698,705,706,710,711,713,714,719,722,724,732,735,737,743,745,746,752,759,761,763,770,772,785,790,799,812,815,816,822,826,829,833,835,840,843,844,849,850,853,856,867,868,876,878,897,903,909,910,923,932,936,942,950,951,952,955,956,965,978,982,1009,1010,1019,1027,1032,1049,1050,1058,1060,1068,1075,1076,1085,1086,1108,1109,1110,1111,1118,1120,1124,1126,1130,1144,1145,1152,1155,1159,1166,1170,1171,1173,1175,1183,1201,1210,1211,1220,1226,1230,1233,1234,1235,1240,1241,1244,1256,1257,1266,1282,1289,1294,1296,1300,1318,1320,1330,1336,1342,1349,1351,1353,1357,1359,1362,1363,1373,1378,1382,1390,1400,1407,1415,1427,1435,1449,1457,1462,1466,1468,1474,1478,1479,1491,1505,1507,1513,1527,1552,1555,1556,1563,1576,1580,1581,1596,1608,1609,1616,1623,1629,1637,1638,1654,1664,1665,1690,1700,1711,1738,1745,1748,1750,1752,1754,1756,1758,1768,1778,1781,1806,1823,1827,1829,1833,1836,1840,1841,1842,1851,1857,1872,1877,1890,1895,1897,1902,1909,1911,1912,1916,1924,1927,1931,1934,1936,1941,1946,1963,1965,1974,1978,1993,1995,2002,2005,2007,2022,2030,2048,2049,2054,2086,2088,2093,2102,2112,2117,2118,2138,2149,2151,2158,2159,2161,2165,2196,2199,2205,2213,2219,2222,2223,2233,2235,2239,2251,2275,2277,2284,2292,2296,2314,2316,2317,2320,2321,2334,2336,2337,2341,2343,2352,2355,2356,2357,2370,2380,2384,2385,2393,2397,2400,2406,2410,2420,2424,2435,2439,2440,2445,2449,2453,2472,2481,2487,2492,2496,2498,2499,2500,2506,2513,2515,2525,2536,2537,2547,2566,2578,2592,2593,2619,2623,2632,2635,2637,2646,2650,2651,2654,2657,2663)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = $false
}
